$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the trailing rows that no longer exist in the updated dataset (37 -> 33 rows)
$ws.Range("A34:B37").ClearContents()

# Update cluster names and active-case counts with the refreshed report data
$ws.Range("A2").Value = "3398 BlueCross Elly Kay Mordialloc"
$ws.Range("B2").Value = 29
$ws.Range("A3").Value = "3601 Baptcare Westhaven community"
$ws.Range("B3").Value = 13
$ws.Range("A4").Value = "3653 Fronditha Thalpori St Albans Aged Care"
$ws.Range("B4").Value = 20
$ws.Range("A5").Value = "3749 Rosebrook - McKenzie Aged Care Rosebud"
$ws.Range("B5").Value = 12
$ws.Range("A6").Value = "3975 Aurrum Aged Care Brunswick"
$ws.Range("B6").Value = 11
$ws.Range("A7").Value = "4257 BlueCross The Gables Camberwell"
$ws.Range("B7").Value = 16
$ws.Range("A8").Value = "4295 Hope Aged Care Sunshine West"
$ws.Range("B8").Value = 14
$ws.Range("A9").Value = "4314 Estia Health Ardeer"
$ws.Range("B9").Value = 13
$ws.Range("A10").Value = "44095 Myrniong Primary School Myrniong"
$ws.Range("B10").Value = 13
$ws.Range("A11").Value = "44404 Castlemaine North Primary School Castlemaine"
$ws.Range("B11").Value = 26
$ws.Range("A12").Value = "44593 Torquay P-6 College Torquay"
$ws.Range("B12").Value = 12
$ws.Range("A13").Value = "44622 Grey Street Primary School Traralgon"
$ws.Range("B13").Value = 14
$ws.Range("A14").Value = "44631 Mount Evelyn Primary School"
$ws.Range("B14").Value = 22
$ws.Range("A15").Value = "44642 Irymple South Primary School Irymple South"
$ws.Range("B15").Value = 20
$ws.Range("A16").Value = "4479 Whittlesea Lodge Whittlesea"
$ws.Range("B16").Value = 16
$ws.Range("A17").Value = "44893 Greenhills Primary School Greensborough"
$ws.Range("B17").Value = 12
$ws.Range("A18").Value = "45168 Ranfurly Primary School Mildura"
$ws.Range("B18").Value = 18
$ws.Range("A19").Value = "45275 Lalor Gardens Primary School Lalor"
$ws.Range("B19").Value = 11
$ws.Range("A20").Value = "45305 Lockington Consolidated School Lockington"
$ws.Range("B20").Value = 16
$ws.Range("A21").Value = "52390 Our Lady of the Way Catholic Primary School Wallan"
$ws.Range("B21").Value = 44
$ws.Range("A22").Value = "52777 Mirripoa Primary School Mount Duneed School Camp"
$ws.Range("B22").Value = 32
$ws.Range("A23").Value = "Confirmed Omicron Sircuit Bar Fitzroy"
$ws.Range("B23").Value = 14
$ws.Range("A24").Value = "Confirmed Omicron Variant The Peel Hotel Collingwood"
$ws.Range("B24").Value = 14
$ws.Range("A25").Value = "Green Gables Lodge Warburton"
$ws.Range("B25").Value = 21
$ws.Range("A26").Value = "JBS Australia Brooklyn"
$ws.Range("B26").Value = 12
$ws.Range("A27").Value = "PGL Camp Rumbug Foster North"
$ws.Range("B27").Value = 10
$ws.Range("A28").Value = "St Brigid's Parish Primary School Mordialloc"
$ws.Range("B28").Value = 12
$ws.Range("A29").Value = "St Pauls Cathedral"
$ws.Range("B29").Value = 42
$ws.Range("A30").Value = "St Vincents Hospital Melbourne Emergency Department Fitzroy"
$ws.Range("B30").Value = 14
$ws.Range("A31").Value = "St. Vincent's Hospital Melbourne Fitzroy"
$ws.Range("B31").Value = 16
$ws.Range("A32").Value = "Sunny Ridge Strawberry Farm Main Ridge"
$ws.Range("B32").Value = 14
$ws.Range("A33").Value = "Warburton Lodge Warburton"
$ws.Range("B33").Value = 14
